$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete old data rows 2 and 3 (shifts remaining rows up by 2)
$ws.Rows("2:3").Delete()

# Append 12 new rows of accelerometer data after the existing data (now ending at row 19)
$newData = @(
    @(1.290312767028809, -2.825088739395142, -12.40687561035156),
    @(-6.366491794586182, -13.87612628936768, -13.16234016418457),
    @(-4.204625129699707, -2.39580774307251, -2.694075107574463),
    @(3.387409210205078, -0.1273889541625976, -3.40644645690918),
    @(-30.62431526184082, -35.95574188232422, -0.9927592277526855),
    @(31.11456108093262, 10.69985771179199, -13.1208086013794),
    @(-20.68763160705566, -28.60056304931641, -20.36888885498047),
    @(27.19162559509277, -76.43211364746094, -1.331388473510742),
    @(2.005526781082153, 2.733113288879395, -6.29358959197998),
    @(-72.39521789550781, -77.27358245849609, 66.16534423828125),
    @(-17.93119621276855, 35.98200607299805, 4.464548587799072),
    @(-0.0441454052925109, -11.34683990478516, -10.79495239257812)
)

$startRow = 20
for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newData[$i][0]
    $ws.Cells.Item($row, 2).Value = $newData[$i][1]
    $ws.Cells.Item($row, 3).Value = $newData[$i][2]
}
